$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column index map:
# A=1 B=2 C=3 D=4 E=5 F=6 G=7 H=8 I=9 J=10 K=11 L=12 M=13 N=14 O=15 P=16
# Q=17 R=18 S=19 T=20 U=21 V=22 W=23 X=24 Y=25 Z=26 AA=27 AB=28 AC=29

# ---------------------------------------------------------------
# Row 236 (id 234): update to new match data
# ---------------------------------------------------------------
$r = 236
$ws.Cells.Item($r, 2).Value = 7013885
$ws.Cells.Item($r, 6).Value = "La Luz"
$ws.Cells.Item($r, 7).Value = "Atletico Fenix Montevideo"
$ws.Cells.Item($r, 8).Value = 0
$ws.Cells.Item($r, 9).Value = 2
$ws.Cells.Item($r, 10).Value = "A"
$ws.Cells.Item($r, 11).Value = 3
$ws.Cells.Item($r, 12).Value = 3
$ws.Cells.Item($r, 13).Value = 2.4
$ws.Cells.Item($r, 14).Value = 2.9
$ws.Cells.Item($r, 15).Value = 2.75
$ws.Cells.Item($r, 16).Value = 2.6
$ws.Cells.Item($r, 17).Value = 0
$ws.Cells.Item($r, 18).Value = 2.025
$ws.Cells.Item($r, 19).Value = 1.825
$ws.Cells.Item($r, 20).Value = 2
$ws.Cells.Item($r, 21).Value = 2.025
$ws.Cells.Item($r, 22).Value = 1.825
$ws.Cells.Item($r, 24).Value = -1
$ws.Cells.Item($r, 25).Value = 1.6
$ws.Cells.Item($r, 27).Value = 0.825
$ws.Cells.Item($r, 28).Value = 0
$ws.Cells.Item($r, 29).Value = -0

# ---------------------------------------------------------------
# Row 239 (id 237): update to new match data
# ---------------------------------------------------------------
$r = 239
$ws.Cells.Item($r, 2).Value = 7013409
$ws.Cells.Item($r, 6).Value = "Nacional De Football"
$ws.Cells.Item($r, 7).Value = "Torque"
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Value = 1
$ws.Cells.Item($r, 10).Value = "D"
$ws.Cells.Item($r, 11).Value = 1.666
$ws.Cells.Item($r, 12).Value = 3.9
$ws.Cells.Item($r, 13).Value = 4.5
$ws.Cells.Item($r, 14).Value = 1.615
$ws.Cells.Item($r, 15).Value = 4
$ws.Cells.Item($r, 16).Value = 4.75
$ws.Cells.Item($r, 17).Value = -0.75
$ws.Cells.Item($r, 18).Value = 1.8
$ws.Cells.Item($r, 19).Value = 2.05
$ws.Cells.Item($r, 20).Value = 2.75
$ws.Cells.Item($r, 21).Value = 1.95
$ws.Cells.Item($r, 22).Value = 1.9
$ws.Cells.Item($r, 24).Value = 3
$ws.Cells.Item($r, 25).Value = -1
$ws.Cells.Item($r, 27).Value = 1.05
$ws.Cells.Item($r, 28).Value = -1
$ws.Cells.Item($r, 29).Value = 0.8999999999999999

# ---------------------------------------------------------------
# Row 249 (id 247): update to new match data (adds FTHG/FTAG/FTR
# and PL_AhOver/PL_AhUnder which previously did not exist)
# ---------------------------------------------------------------
$r = 249
$ws.Cells.Item($r, 2).Value = 7825145
$ws.Cells.Item($r, 5).Value = 45345.85416666666
$ws.Cells.Item($r, 6).Value = "Montevideo Wanderers"
$ws.Cells.Item($r, 7).Value = "Racing Club de Montevideo"
$ws.Cells.Item($r, 8).Value = 0
$ws.Cells.Item($r, 9).Value = 1
$ws.Cells.Item($r, 10).Value = "A"
$ws.Cells.Item($r, 11).Value = 2.5
$ws.Cells.Item($r, 12).Value = 3.2
$ws.Cells.Item($r, 13).Value = 2.75
$ws.Cells.Item($r, 14).Value = 2.6
$ws.Cells.Item($r, 15).Value = 3.2
$ws.Cells.Item($r, 16).Value = 2.625
$ws.Cells.Item($r, 17).Value = 0
$ws.Cells.Item($r, 18).Value = 1.95
$ws.Cells.Item($r, 19).Value = 1.9
$ws.Cells.Item($r, 20).Value = 2
$ws.Cells.Item($r, 21).Value = 1.85
$ws.Cells.Item($r, 22).Value = 2
$ws.Cells.Item($r, 23).Value = -1
$ws.Cells.Item($r, 24).Value = -1
$ws.Cells.Item($r, 25).Value = 1.625
$ws.Cells.Item($r, 26).Value = -1
$ws.Cells.Item($r, 27).Value = 0.8999999999999999
$ws.Cells.Item($r, 28).Value = -1
$ws.Cells.Item($r, 29).Value = 1

# ---------------------------------------------------------------
# Row 250 (id 248): update to new match data (adds FTHG/FTAG/FTR
# and PL_AhOver/PL_AhUnder which previously did not exist)
# ---------------------------------------------------------------
$r = 250
$ws.Cells.Item($r, 2).Value = 7825143
$ws.Cells.Item($r, 5).Value = 45346.70833333334
$ws.Cells.Item($r, 6).Value = "Boston River"
$ws.Cells.Item($r, 7).Value = "Danubio"
$ws.Cells.Item($r, 8).Value = 0
$ws.Cells.Item($r, 9).Value = 0
$ws.Cells.Item($r, 10).Value = "D"
$ws.Cells.Item($r, 11).Value = 2.75
$ws.Cells.Item($r, 12).Value = 3.1
$ws.Cells.Item($r, 13).Value = 2.5
$ws.Cells.Item($r, 14).Value = 2.55
$ws.Cells.Item($r, 15).Value = 3.1
$ws.Cells.Item($r, 16).Value = 2.7
$ws.Cells.Item($r, 17).Value = 0
$ws.Cells.Item($r, 20).Value = 2.25
$ws.Cells.Item($r, 23).Value = -1
$ws.Cells.Item($r, 24).Value = 2.1
$ws.Cells.Item($r, 25).Value = -1
$ws.Cells.Item($r, 27).Value = -0
$ws.Cells.Item($r, 28).Value = -1
$ws.Cells.Item($r, 29).Value = 0.875
